$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-01 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-02 Saturday", 2)

$d.Content.Find.Execute("458÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "386÷3=", 2)
$d.Content.Find.Execute("141÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "451÷8=", 2)
$d.Content.Find.Execute("137÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "660÷6=", 2)
$d.Content.Find.Execute("329÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "446÷3=", 2)
$d.Content.Find.Execute("764÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "133÷6=", 2)
$d.Content.Find.Execute("750÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "455÷4=", 2)
$d.Content.Find.Execute("611÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "343÷2=", 2)
$d.Content.Find.Execute("845÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "392÷4=", 2)
$d.Content.Find.Execute("722÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "695÷9=", 2)
$d.Content.Find.Execute("222÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "810÷4=", 2)
$d.Content.Find.Execute("995÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "303÷7=", 2)
$d.Content.Find.Execute("999÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "612÷8=", 2)
$d.Content.Find.Execute("370÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "960÷9=", 2)
$d.Content.Find.Execute("884÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "743÷5=", 2)
$d.Content.Find.Execute("236÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "409÷5=", 2)
$d.Content.Find.Execute("782÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "749÷3=", 2)
$d.Content.Find.Execute("524÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "943÷8=", 2)
$d.Content.Find.Execute("809÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "371÷5=", 2)
$d.Content.Find.Execute("276÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "265÷7=", 2)
$d.Content.Find.Execute("514÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "519÷2=", 2)
$d.Content.Find.Execute("135÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "875÷5=", 2)
$d.Content.Find.Execute("314÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "100÷9=", 2)
$d.Content.Find.Execute("451÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "356÷5=", 2)
$d.Content.Find.Execute("608÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "477÷6=", 2)
$d.Content.Find.Execute("319÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "296÷5=", 2)
